$wb = $excel.ActiveWorkbook

# --- Map sheet: currently the selected/active tab; Flags becomes active instead ---
$wsFlags = $wb.Worksheets.Item("Flags")

# --- Flags sheet: new header row (row 7) and data rows (8-15) ---
# Values are written in this precise order so that new shared strings are
# interned in the same sequence as the target workbook (B7, D7, C7, A7, F7, G7, E7, then H11-H13).
$wsFlags.Range("B7").Value = "Prepared 3rd"
$wsFlags.Range("D7").Value = "Direction change 2"
$wsFlags.Range("C7").Value = "Direction change 1"
$wsFlags.Range("A7").Value = "Leap2leap"
$wsFlags.Range("F7").Value = "Both 3rds"
$wsFlags.Range("G7").Value = "Prepared 2 3rds"
$wsFlags.Range("E7").Value = "Leap <6"

$wsFlags.Range("A8").Value = 1
$wsFlags.Range("C8").Value = 0
$wsFlags.Range("F8").Value = 1
$wsFlags.Range("G8").Value = 1
$wsFlags.Range("H8").Value = "Two 3rds after 6/8"

$wsFlags.Range("A9").Value = 1
$wsFlags.Range("C9").Value = 0
$wsFlags.Range("F9").Value = 1
$wsFlags.Range("G9").Value = 0
$wsFlags.Range("H9").Value = "Two 3rds"

$wsFlags.Range("F10").Value = 0
$wsFlags.Range("H10").Value = "Leap chain"

$wsFlags.Range("A11").Value = 0
$wsFlags.Range("B11").Value = 1
$wsFlags.Range("C11").Value = 0
$wsFlags.Range("H11").Value = "Unresolved prepared 3rd"

$wsFlags.Range("C12").Value = 0
$wsFlags.Range("D12").Value = 1
$wsFlags.Range("E12").Value = 1
$wsFlags.Range("H12").Value = "Late leap resolve <6"

$wsFlags.Range("C13").Value = 0
$wsFlags.Range("D13").Value = 1
$wsFlags.Range("E13").Value = 0
$wsFlags.Range("H13").Value = "Late leap resolve >5"

$wsFlags.Range("C14").Value = 0
$wsFlags.Range("D14").Value = 0
$wsFlags.Range("H14").Value = "Leap unresolved"

$wsFlags.Range("H15").Value = "Leap2leap"

# Column G needs to accommodate the new "Prepared 2 3rds" header.
$wsFlags.Columns.Item(7).AutoFit()

# Page setup (print settings) for the Flags sheet.
$wsFlags.PageSetup.PaperSize = 9
$wsFlags.PageSetup.Orientation = 1

# --- Selection / active-sheet bookkeeping ---
# Map sheet (previously tabSelected) loses the flag; Flags sheet (previously
# plain) becomes the selected / active tab.
[void]$wsFlags.Range("H16").Select()
[void]$wsFlags.Activate()
